$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value to insert into column B (quarter Q0 / most-recent error)
$newB = @{
    2 = 2.267328492325126
    3 = 9.885106142864725
    4 = -9.085154484867306
    5 = -0.02759857525660239
    6 = 0.6014450472570072
    7 = -1.308035509750171
    8 = -2.058868485289545
    9 = 0.8913479099652445
    10 = -0.7322633397437844
    11 = 0.2703549766394939
    12 = -1.355327161308811
    13 = 1.651602845777944
    14 = 0.3282974736644749
    15 = 0.7356582956163805
    16 = 0.1181882633125878
    17 = 0.7543890506736601
    18 = -0.1543252035281459
    19 = 0.2293445564577608
    20 = 0.2201546830999171
    21 = 0.314534851581486
    22 = -0.5970339283829468
    23 = 0.1550649743121164
    24 = -0.1624199859130616
}

# For each data row (2..24), shift existing values in columns B:J (2..10) right into C:K (3..11),
# working from the rightmost column down to B so we don't clobber data before reading it.
# This mirrors prepending a new quarter's error estimate (Q0) and dropping the oldest tracked quarter
# (old column K / Q9) to keep the fixed 10-quarter window.
for ($r = 2; $r -le 24; $r++) {
    for ($c = 10; $c -ge 2; $c--) {
        $srcVal = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r, $c + 1).Value = $srcVal
    }
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
